$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that should be bumped from
# 45171 (2023-09-02) to 45172 (2023-09-03) for every data row (rows 2-387).
$lastRow = 387
$ws.Range("C2:C$lastRow").Value = 45172
